$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# fix(Klein): use EMK for c*phi_E in Leerlaufversuch
# E = EMK * 60 / n  (column B * 60 / column A)
$ws.Range("E2").Value = 0.1332
$ws.Range("E3").Value = 0.132
$ws.Range("E4").Value = 0.132
$ws.Range("E5").Value = 0.1329
$ws.Range("E6").Value = 0.132
